# Insert two new rows of feature-comparison data into Sheet1.
# Row order before edit (rows 7-26):
#   7  franzosa_ControlvsCD_Age
#   8  franzosa_ControlvsCD_Fp
#   9  franzosa_ControlvsDisease_Age
#   10 franzosa_ControlvsDisease_ConvDisease
#   11 franzosa_ControlvsDisease_Fp
#   12 franzosa_ControlvsUC_Age
#   13 franzosa_ControlvsUC_Fp
#   14 yachida_age
#   ...
# After the edit, two brand-new rows are inserted (rest shift down by one for
# each insertion, retaining their original values):
#   - "franzosa_ControlvsCD_ConvCD" right after "franzosa_ControlvsCD_Age"
#     (becomes new row 8; old row 8 "franzosa_ControlvsCD_Fp" becomes row 9)
#   - "franzosa_ControlvsUC_ConvUC" right after "franzosa_ControlvsUC_Age"
#     (becomes new row 14; old row 13 "franzosa_ControlvsUC_Fp" becomes row 15)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "franzosa_ControlvsCD_ConvCD" as a new row 8 ---
$ws.Rows.Item(8).Insert()
$ws.Range("A8").Value = "franzosa_ControlvsCD_ConvCD"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.4
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6
$ws.Range("H8").Value = 0.6

# --- Insert "franzosa_ControlvsUC_ConvUC" as a new row 14 ---
# (franzosa_ControlvsUC_Age is now row 13 after the first insertion above)
$ws.Rows.Item(14).Insert()
$ws.Range("A14").Value = "franzosa_ControlvsUC_ConvUC"
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0.4
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.6
$ws.Range("H14").Value = 0.6
